$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 need the same style as the other header cells (H1),
# so copy H1's formatting across before writing the new header text.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-7 for new columns I (I0) and J (IF)
$data = @(
    @(6, 9),
    @(7, 9),
    @(10, 10),
    @(8, 9),
    @(1, 3),
    @(1, 2)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
